$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the query text in cell B4: remove the "Not specified in data" default
# fallback for Participant ID and Sample ID, replacing it with an empty string.
$oldText = $ws.Range("B4").Value2
$newText = $oldText.Replace("coalesce(p.participant_id, 'Not specified in data')", "coalesce(p.participant_id, '')")
$newText = $newText.Replace("coalesce(samp.sample_id, 'Not specified in data')", "coalesce(samp.sample_id, '')")
$ws.Range("B4").Value2 = $newText

# Update the sheet view: clear the scrolled topLeftCell and move the
# active selection from D4 to B2.
$ws.Activate()
$ws.Range("B2").Select()
